$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cell format for columns that hold text-formatted numbers
# (prices/links/coin names are stored as plain text strings, not numeric values)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.543.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.91%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4636"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.04%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.63%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.85"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.24%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07885"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9960"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.46"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.58%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.99%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.950"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.124"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.42%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.23%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.96%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.47%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.541.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.382"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.30%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.069.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.81%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.111"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.405"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.92%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9751"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09394"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.85%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.293"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.338"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06032"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.35%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.303"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5881"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1865"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.33"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.234"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5587"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.95%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06696"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.90"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.35%  "
